# Auto update Excel log
# Appends newly-logged sensor/alert rows to the bottom of five sheets
# (ALERTS, PIR, Humidity, Temperature, Proximity) for the 2026-02-01 18:00 hour.
$wb = $excel.ActiveWorkbook

# --- ALERTS: append rows 16-18 ---
$ws = $wb.Worksheets.Item("ALERTS")
$newRows = @(
    ,("'2026-02-01","'18:35:44","'18:00","'Bathroom","'MINIMAL","'MINIMAL ALERT: Bathroom occupied, no motion > 20s.")
    ,("'2026-02-01","'18:35:49","'18:00","'Bathroom","'WARNING","'Bathroom Humidity > 90.0% for 22s with NO MOTION. Alerting.")
    ,("'2026-02-01","'18:36:14","'18:00","'Bathroom","'MINIMAL","'MINIMAL ALERT: Bathroom occupied, no motion > 20s.")
)
$startRow = 16
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c-1]
    }
}

# --- PIR: append rows 81-93 ---
$ws = $wb.Worksheets.Item("PIR")
$newRows = @(
    ,("'2026-02-01","'18:35:17","'18:00","'Bathroom","'No Motion","'Inactive")
    ,("'2026-02-01","'18:35:20","'18:00","'Bathroom","'No Motion","'Inactive")
    ,("'2026-02-01","'18:35:23","'18:00","'Bathroom","'No Motion","'Inactive")
    ,("'2026-02-01","'18:35:28","'18:00","'Bathroom","'No Motion","'Inactive")
    ,("'2026-02-01","'18:35:33","'18:00","'Bathroom","'No Motion","'Inactive")
    ,("'2026-02-01","'18:35:38","'18:00","'Bathroom","'No Motion","'Inactive")
    ,("'2026-02-01","'18:35:43","'18:00","'Bathroom","'No Motion","'Inactive")
    ,("'2026-02-01","'18:35:48","'18:00","'Bathroom","'No Motion","'Inactive")
    ,("'2026-02-01","'18:35:51","'18:00","'Bathroom","'Motion Detected","'Active")
    ,("'2026-02-01","'18:35:58","'18:00","'Bathroom","'No Motion","'Inactive")
    ,("'2026-02-01","'18:36:03","'18:00","'Bathroom","'No Motion","'Inactive")
    ,("'2026-02-01","'18:36:08","'18:00","'Bathroom","'No Motion","'Inactive")
    ,("'2026-02-01","'18:36:13","'18:00","'Bathroom","'No Motion","'Inactive")
)
$startRow = 81
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c-1]
    }
}

# --- Humidity: append rows 139-148 ---
$ws = $wb.Worksheets.Item("Humidity")
$newRows = @(
    ,("'2026-02-01","'18:35:19","'18:00","'Bathroom","'99.9%","'Active")
    ,("'2026-02-01","'18:35:20","'18:00","'Bathroom","'99.9%","'Active")
    ,("'2026-02-01","'18:35:24","'18:00","'Bathroom","'99.9%","'Active")
    ,("'2026-02-01","'18:35:29","'18:00","'Bathroom","'99.9%","'Active")
    ,("'2026-02-01","'18:35:34","'18:00","'Bathroom","'99.9%","'Active")
    ,("'2026-02-01","'18:35:39","'18:00","'Bathroom","'98.5%","'Active")
    ,("'2026-02-01","'18:35:45","'18:00","'Bathroom","'99.9%","'Active")
    ,("'2026-02-01","'18:35:50","'18:00","'Bathroom","'99.9%","'Active")
    ,("'2026-02-01","'18:35:55","'18:00","'Bathroom","'99.9%","'Active")
    ,("'2026-02-01","'18:35:59","'18:00","'Bathroom","'98.6%","'Active")
)
$startRow = 139
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c-1]
    }
}

# --- Temperature: append rows 139-148 ---
$ws = $wb.Worksheets.Item("Temperature")
$newRows = @(
    ,("'2026-02-01","'18:35:19","'18:00","'Bathroom","'29.9C","'Active")
    ,("'2026-02-01","'18:35:21","'18:00","'Bathroom","'29.9C","'Active")
    ,("'2026-02-01","'18:35:25","'18:00","'Bathroom","'29.8C","'Active")
    ,("'2026-02-01","'18:35:30","'18:00","'Bathroom","'29.8C","'Active")
    ,("'2026-02-01","'18:35:35","'18:00","'Bathroom","'30.0C","'Active")
    ,("'2026-02-01","'18:35:40","'18:00","'Bathroom","'30.1C","'Active")
    ,("'2026-02-01","'18:35:45","'18:00","'Bathroom","'30.5C","'Active")
    ,("'2026-02-01","'18:35:50","'18:00","'Bathroom","'30.3C","'Active")
    ,("'2026-02-01","'18:35:55","'18:00","'Bathroom","'30.5C","'Active")
    ,("'2026-02-01","'18:36:00","'18:00","'Bathroom","'30.5C","'Active")
)
$startRow = 139
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c-1]
    }
}

# --- Proximity: append rows 50-52 ---
$ws = $wb.Worksheets.Item("Proximity")
$newRows = @(
    ,("'2026-02-01","'18:35:17","'18:00","'Bathroom Door","'ENTER","'User ENTERED Bathroom")
    ,("'2026-02-01","'18:35:18","'18:00","'Bathroom Door","'EXIT","'User EXITED Bathroom")
    ,("'2026-02-01","'18:35:22","'18:00","'Bathroom Door","'ENTER","'User ENTERED Bathroom")
)
$startRow = 50
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c-1]
    }
}

